# CIV-11205 updated templates with court full address
#
# The "This order is made by <<judgeNameTitle>> on <<submittedOn>> at
# <<courtName>>." paragraph now spells out the court's full address
# instead of just its name, i.e. "<<courtName>>" becomes
# "<<siteName>> - <<address>> - <<postcode>>".

$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "<<courtName>>.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "<<siteName>> - <<address>> - <<postcode>>.",
    2
)

Write-Output "Replaced courtName placeholder: $found"
